# Update the quote-number cells on the two "QuoteSave" sheets.
# A new batch of randomly generated quote numbers was produced by the
# application; the most recently generated code for each sheet is written
# into that sheet's A1 cell (all other cell contents are left untouched).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("BrokerPCHQuoteNo")
$ws1.Range("A1").Value = "ZHE09143"

$ws3 = $wb.Worksheets.Item("BrokerBCHQuoteNo")
$ws3.Range("A1").Value = "XOR36155"
